$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Adjust column widths: D, E, F
# (Excel's ColumnWidth property reads/writes ~0.83 chars lower than the
# raw OOXML <col width> value for this sheet's font metrics, so offset
# the target widths by +0.83 to land exactly on 11 / 22 / 18.)
$ws.Columns.Item(4).ColumnWidth = 10.17
$ws.Columns.Item(5).ColumnWidth = 21.17
$ws.Columns.Item(6).ColumnWidth = 17.17

# Row 2 (240X120 PORCELANATO)
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1041.16
$ws.Range("F2").Value = 0

# Row 3 (240X80 PORCELANATO)
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 8668.91
$ws.Range("F3").Value = 0

# Row 4 (FREGADEROS DE COCINA)
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 372.993863046034
$ws.Range("F4").Value = 0

# Row 6 (GRIFERIAS)
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 106.82
$ws.Range("F6").Value = 0

# Row 8 (LAVABOS)
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 625
$ws.Range("F8").Value = 0

# Row 12 (PANELES DECORATIVOS)
$ws.Range("C12").Value = 100
$ws.Range("E12").Value = 100

# Row 13 (PANELES PU)
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 20
$ws.Range("F13").Value = 0

# Row 14 (PANELES PVC)
$ws.Range("C14").Value = 100
$ws.Range("E14").Value = 100

# Row 15 (PIEDRA SINTERIZADA)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 2501.01
$ws.Range("F15").Value = 0

# Row 16 (PORCELANATO)
$ws.Range("C16").Value = 36056.7
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 36056.7
$ws.Range("F16").Value = 0

# Row 19 (TOTAL)
$ws.Range("C19").Value = 55023.16386304603
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 55023.16386304603
$ws.Range("F19").Value = 0
